# Fill in the missing assignment scores (A2, A3, A4 columns = H, I, J)
# for several students whose rows previously had blank cells.
# Totals (K), Calificación (L) and Puntos (M) are formula-driven and
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1

# Row 13
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0

# Row 15
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1

# Row 17 (only J was missing)
$ws.Range("J17").Value = 1

# Row 18 (H and I were missing)
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1

# Row 22
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 0

# Update the active selection to match the saved view state (J1)
$ws.Range("J1").Select() | Out-Null
